# Analyse de la partie MenuPrincipale
# Fill in the journal entries for the day of 2022-05-03 (serial date 44684)
# in rows 7 through 12 of the "Journal" table on Feuil1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = -4122

# --- Row 7 : 08:00 -> 08:17 ---
$ws.Range("A7").Value = 44684
$ws.Range("B7").Value = 0.33333333333333331
$ws.Range("C7").Value = 0.34513888888888888
$ws.Range("E7").Value = "Srprint Mise en Forme"
$ws.Range("F7").Value = "Envoie d'une proposition de mise en forme au chef de projet"

# --- Row 8 : 08:17 -> 08:30 ---
$ws.Range("A8").Value = 44684
$ws.Range("B8").Value = 0.34513888888888888
$ws.Range("C8").Value = 0.35416666666666669
$ws.Range("E8").Value = "Remise en forme`ndu planning"
$ws.Range("F8").Value = "Réception du rendez-vous de l'expert 2 demandant d'ajuster la planification en conséquence"

# --- Row 9 : 08:30 -> 10:23 ---
$ws.Range("A9").Value = 44684
$ws.Range("B9").Value = 0.35416666666666669
$ws.Range("C9").Value = 0.43263888888888885
$ws.Range("E9").Value = "Documentation"
$ws.Range("F9").Value = "Remplissage de la partie Analyse préliminaire`npartie beaucoup plus longue et laborieuse qu'envisagée de par les terme a mettre dans le glossaire et les acronyme à mettre en bas de page"
$ws.Range("G9").Value = "le week end sera LONG!"

# --- Row 10 : 10:23 -> 11:30 ---
$ws.Range("A10").Value = 44684
$ws.Range("B10").Value = 0.43263888888888885
$ws.Range("C10").Value = 0.47916666666666669
$ws.Range("E10").Value = "Analyse Menu Principale"
$ws.Range("F10").Value = "Analyse prélinimaire"
$ws.Range("G10").Value = "ötude de l'interface du jeu New World"

# --- Row 11 : 11:30 -> 12:15 ---
$ws.Range("A11").Value = 44684
$ws.Range("B11").Value = 0.47916666666666669
$ws.Range("C11").Value = 0.51041666666666663
$ws.Range("E11").Value = "Entretien`n avec le Chef de projet"
$ws.Range("F11").Value = "Question sur l'organistion des sprint`nConvention de nommage`nAnalyse et structure du Canvas"
$ws.Range("G11").Value = "Utilisation de Icescrum`nChaque élément du code doit être fixé`nNe pas aller trop loin et ne mettre que des images ciblant précisément le sujet qui support l'illustration`n"

# --- Row 12 : 13:30 -> (ongoing, no end time yet) ---
$ws.Range("A12").Value = 44684
$ws.Range("B12").Value = 0.5625
$ws.Range("E12").Value = "Convention de Nommage"

# A9/A11 previously carried a stray "time" number format (leftover from copy/paste);
# now that they hold real dates, re-apply the same date formatting used by the
# rest of column A (e.g. A2). Likewise E8/E11 pick up the sibling "Description"
# formatting already used by E2 elsewhere in the column.
$ws.Range("A2").Copy()
$ws.Range("A9").PasteSpecial($xlPasteFormats)
$ws.Range("A11").PasteSpecial($xlPasteFormats)

$ws.Range("E2").Copy()
$ws.Range("E8").PasteSpecial($xlPasteFormats)
$ws.Range("E11").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# Leave the cursor where the author left off after typing the new entries.
[void]$ws.Range("C12").Select()
